$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update category text for rows that previously selected other options,
# now reflecting the stricter Coding threshold re-run (they no longer
# qualify, so they did not choose "Brainstorming & personal ideas / fun").
$ws.Range("B2").Value = "I did not choose " + [char]8220 + "Brainstorming & personal ideas / fun" + [char]8221
$ws.Range("B4").Value = "I did not choose " + [char]8220 + "Brainstorming & personal ideas / fun" + [char]8221

# Update timestamps for all data rows (2-9) to reflect the re-run time.
$newTimestamp = 45841.64456272784
$ws.Range("D2:D9").Value = $newTimestamp
